$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A91").Value = "Globo"
$ws.Range("B91").Value = "RJ TV 2"
$ws.Range("C91").Value = "Governo"
$ws.Range("D91").Value = "2025-04-08T19:31"
$ws.Range("E91").Value = "Neutro"
$ws.Range("F91").Value = "Reforma administrativa suspensa na Justiça. Decisão suspende tramitação na Câmara de reforma proposta pela Prefeitura. Repórter *ao vivo* em frente à Câmara Municipal. Repórter lembrou que, no início do ano, conversou com o presidente da Casa, Fred Rangel, que falou que esse projeto era uma das principais propostas a serem votadas. Projeto pode causar acréscimo de 10% na folha de pagamento. Entrevista com a vereadora Thamires Rangel, autora do pedido de suspensão, e com o líder do governo na Câmara, vereador Juninho Virgílio. Previsão é de que a reforma fosse para o plenário hoje ou amanhã. "

$ws.Range("A92").Value = "Record"
$ws.Range("B92").Value = "RJ Record"
$ws.Range("C92").Value = "Governo"
$ws.Range("D92").Value = "2025-04-08T18:26"
$ws.Range("E92").Value = "Neutro"
$ws.Range("F92").Value = "Justiça suspende reforma administrativa da prefeitura que tramitava na Câmara. Repórter *ao vivo*. Decisão é do juiz da 4ª Vara Cível de Campos. Mandado de segurança foi impetrado pela vereadora Thamires Rangel. Na liminar, juiz informa previsão de aumento da despesa mensal de 9,6% com folha de pagamento, que já é de R$ 1 bi por ano. Secretário interino da Transparência não deixa claro a fonte para cobrir essas despesas. Juiz deu prazo de 10 dias para explicar detalhes do que falta no projeto. "
